$wb = $excel.ActiveWorkbook

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 406.81818
$ws.Range("I33").Value = 436.5
$ws.Range("J33").Value = 110
$ws.Range("K33").Value = 436.5
$ws.Range("L33").Value = 110
$ws.Range("M33").Value = -207.5
$ws.Range("N33").Value = -568

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3091.1667
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 2886.75
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 2886.75
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -9394.75

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2133.3333
$ws.Range("I125").Value = 1866.6666
$ws.Range("K125").Value = 16799.9994
$ws.Range("M125").Value = -14339.9994

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1271.4286
$ws.Range("I45").Value = 1250
$ws.Range("K45").Value = 1250
$ws.Range("M45").Value = -873

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 109825.055
$ws.Range("I63").Value = 124601.3
$ws.Range("J63").Value = 1018.1818
$ws.Range("K63").Value = 124601.3
$ws.Range("L63").Value = 1018.1818
$ws.Range("M63").Value = -123915.3
$ws.Range("N63").Value = -2390.1818

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 109825.055
$ws.Range("I66").Value = 124601.3
$ws.Range("J66").Value = 1018.1818
$ws.Range("K66").Value = 623006.5
$ws.Range("L66").Value = 5090.909
$ws.Range("M66").Value = -619574.5
$ws.Range("N66").Value = -11954.909

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1647.05
$ws.Range("I74").Value = 1710.0714
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1710.0714
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -836.0714
$ws.Range("N74").Value = -3248

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1647.05
$ws.Range("I77").Value = 1710.0714
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 8550.357
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -4182.357
$ws.Range("N77").Value = -16236

# BSM!row106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 29057
$ws.Range("J106").Value = 29057
$ws.Range("L106").Value = 29057
$ws.Range("N106").Value = -31581

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6453910
$ws.Range("I31").Value = 2373.5334
$ws.Range("J31").Value = 200000000
$ws.Range("K31").Value = 2373.5334
$ws.Range("L31").Value = 200000000
$ws.Range("M31").Value = -2078.5334
$ws.Range("N31").Value = -200000590

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6453910
$ws.Range("I34").Value = 2373.5334
$ws.Range("J34").Value = 200000000
$ws.Range("K34").Value = 2373.5334
$ws.Range("L34").Value = 200000000
$ws.Range("M34").Value = -2171.5334
$ws.Range("N34").Value = -200000404

# CRP!row69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 11000
$ws.Range("I69").Value = 6500
$ws.Range("K69").Value = 6500
$ws.Range("M69").Value = -5751

# CRP!row72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 11000
$ws.Range("I72").Value = 6500
$ws.Range("K72").Value = 19500
$ws.Range("M72").Value = -15756

# CRP!row93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 17250
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 29500
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 29500
$ws.Range("M93").Value = -3128
$ws.Range("N93").Value = -33244

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2666.1667
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2999.25
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2999.25
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5995.25

# CRP!row103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 13782.857
$ws.Range("I103").Value = 12746.667
$ws.Range("J103").Value = 20000
$ws.Range("K103").Value = 12746.667
$ws.Range("L103").Value = 20000
$ws.Range("M103").Value = -11574.667
$ws.Range("N103").Value = -22344

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 856.2
$ws.Range("I122").Value = 813.6667
$ws.Range("J122").Value = 920
$ws.Range("K122").Value = 2441.0001
$ws.Range("L122").Value = 2760
$ws.Range("M122").Value = 8.999899999999798
$ws.Range("N122").Value = -7660

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2666.1667
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2999.25
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8997.75
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13937.75

# CUL!row4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 83.2
$ws.Range("I4").Value = 83.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 249.6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -137.6
$ws.Range("N4").ClearContents()

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 485.57574
$ws.Range("I97").Value = 310.16666
$ws.Range("J97").Value = 696.06665
$ws.Range("K97").Value = 310.16666
$ws.Range("L97").Value = 696.06665
$ws.Range("M97").Value = 185.83334
$ws.Range("N97").Value = -1688.06665

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1754.3529
$ws.Range("I102").Value = 1730.9
$ws.Range("K102").Value = 1730.9
$ws.Range("M102").Value = -108.9000000000001

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 144899.78
$ws.Range("I132").Value = 183336.19
$ws.Range("K132").Value = 550008.5700000001
$ws.Range("M132").Value = -547478.5700000001

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3742.5
$ws.Range("I7").Value = 3091
$ws.Range("K7").Value = 3091
$ws.Range("M7").Value = -2979

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3742.5
$ws.Range("I126").Value = 3091
$ws.Range("K126").Value = 9273
$ws.Range("M126").Value = -6803

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11470
$ws.Range("I136").Value = 17116.666
$ws.Range("K136").Value = 51349.99800000001
$ws.Range("M136").Value = -48799.99800000001

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1298.2858
$ws.Range("I122").Value = 1436.8462
$ws.Range("J122").Value = 1073.125
$ws.Range("K122").Value = 4310.5386
$ws.Range("L122").Value = 3219.375
$ws.Range("M122").Value = -1860.5386
$ws.Range("N122").Value = -8119.375

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7707.353
$ws.Range("I136").Value = 9166.071
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 27498.213
$ws.Range("M136").Value = -24948.213
$ws.Range("N136").Value = -7800
